$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Profit" column header to "Revenue"
$ws.Range("C1").Value = "Revenue"

# Select cell C2, matching the saved selection state in the target file
$ws.Range("C2").Select()
